# Generate Report for Handback
#
# This mirrors a "regenerate handback report" run: the handoff/handback
# timestamps for the first file (0ef44f93-...) get refreshed for both the
# zh-cn and de-de language sheets, and the Overview sheet's
# "Latest HO Xliff Generate Date" column is recomputed from the newest of
# those per-language handoff timestamps.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")
$overview = $wb.Worksheets.Item("Overview")

# zh-cn sheet: row 2 is the 0ef44f93-... file.
# H = Correspond Handoff Datetime, K = Correspond Handback DateTime
$zhcn.Range("H2").Value = "2016-08-14 17:04:47"
$zhcn.Range("K2").Value = "2016-08-14 17:05:15"

# de-de sheet: row 2 is the 0ef44f93-... file.
$dede.Range("H2").Value = "2016-08-14 17:04:56"
$dede.Range("K2").Value = "2016-08-14 17:05:26"

# Overview sheet: "Latest HO Xliff Generate Date" (column G) reflects the
# newest handoff datetime across all languages for that file.
$overview.Range("G2").Value = "2016-08-14 17:04:56"
